$d = $word.ActiveDocument

$pairs = @(
    @{old = "671×8=5368"; new = "840×7=5880"},
    @{old = "688×2=1376"; new = "993×6=5958"},
    @{old = "395×4=1580"; new = "110×5=550"},
    @{old = "180×7=1260"; new = "488×7=3416"},
    @{old = "191×7=1337"; new = "815×6=4890"},
    @{old = "103×6=618"; new = "215×9=1935"},
    @{old = "194×4=776"; new = "576×3=1728"},
    @{old = "956×9=8604"; new = "555×5=2775"},
    @{old = "330×2=660"; new = "406×4=1624"},
    @{old = "703×8=5624"; new = "254×6=1524"},
    @{old = "981×5=4905"; new = "402×4=1608"},
    @{old = "441×2=882"; new = "846×8=6768"},
    @{old = "725×5=3625"; new = "687×6=4122"},
    @{old = "723×9=6507"; new = "671×7=4697"},
    @{old = "615×3=1845"; new = "621×5=3105"},
    @{old = "925×3=2775"; new = "838×8=6704"},
    @{old = "155×7=1085"; new = "936×7=6552"},
    @{old = "201×7=1407"; new = "682×4=2728"},
    @{old = "444×5=2220"; new = "987×4=3948"},
    @{old = "259×9=2331"; new = "757×8=6056"},
    @{old = "414×6=2484"; new = "974×8=7792"},
    @{old = "670×6=4020"; new = "438×7=3066"},
    @{old = "607×7=4249"; new = "331×7=2317"},
    @{old = "230×4=920"; new = "850×5=4250"},
    @{old = "568×7=3976"; new = "665×9=5985"}
)

foreach ($pair in $pairs) {
    $d.Content.Find.Execute($pair.old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $pair.new, 2)
}
